# Custody Query Results Mapping - add Driver License ID / Driver License Source / FBI ID rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Insert 3 new blank rows above row 16 (old row 16 "First name" etc. shifts down to row 19).
# Inserting here copies formatting (styles) from the row above (row 15), matching the
# target styles (A=19, B=20, C=15, D=15, E=15, F=4).
$ws.Rows.Item(16).Resize(3).Insert()

# The insert also copied an (empty) F column cell from row 15's formatting; the target
# rows 16-18 only have cells A-E, so drop the spurious F cells.
$ws.Range("F16:F18").Clear()

# Row 16: Driver's License ID
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = "Driver's License ID"
$ws.Range("C16").Value = "Driver License ID"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:DriverLicense/j:DriverLicenseCardIdentification/nc:IdentificationID"

# Row 17: Driver License Source
$ws.Range("B17").Value = "Driver License Source"
$ws.Range("C17").Value = "Driver License Source"
$ws.Range("C17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:DriverLicense/j:DriverLicenseCardIdentification/nc:IdentificationSourceText"

# Row 18: FBI ID
$ws.Range("B18").Value = "FBI ID"
$ws.Range("C18").Value = "FBI ID"
$ws.Range("E18").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/nc:Person[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/j:PersonAugmentation/j:PersonFBIIdentification/nc:IdentificationID"

# Match the row height (56, same as the neighboring wrapped-text rows) on the new rows.
$ws.Range("A16:A18").RowHeight = 56

# Update the frozen-pane scroll position and the active cell selection to match the
# post-edit view (user scrolled up and selected C17 while reviewing the new rows).
$ws.Range("C17").Select()

# Move the saved window position, matching the author's on-screen window placement.
$win = $excel.ActiveWindow
$win.Left = 3480
$win.Top = 3220
